$d = $word.ActiveDocument

# Mapping of old text -> new text (unique across the document, so a
# straightforward Find/Replace for each pair is unambiguous).
$replacements = [ordered]@{
    "2025-12-14 Sunday" = "2025-12-15 Monday"
    "192÷7="            = "136÷2="
    "263÷4="            = "537÷5="
    "951÷9="            = "154÷7="
    "519÷8="            = "905÷3="
    "486÷4="            = "204÷2="
    "963÷3="            = "225÷8="
    "780÷5="            = "110÷7="
    "571÷6="            = "151÷4="
    "521÷4="            = "717÷9="
    "860÷5="            = "560÷5="
    "807÷2="            = "222÷9="
    "784÷8="            = "120÷8="
    "836÷3="            = "308÷8="
    "602÷6="            = "755÷9="
    "164÷4="            = "201÷6="
    "819÷4="            = "396÷3="
    "618÷3="            = "430÷9="
    "748÷8="            = "473÷3="
    "481÷6="            = "257÷3="
    "502÷2="            = "253÷4="
    "378÷5="            = "165÷5="
    "437÷3="            = "142÷7="
    "932÷7="            = "835÷3="
    "837÷4="            = "230÷8="
    "529÷4="            = "276÷7="
}

foreach ($old in $replacements.Keys) {
    $new = $replacements[$old]
    $range = $d.Content
    $range.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2)
}
